# Applies the "Renaissance ZGC chi-square heap-2G" docx edits described by
# the commit: rewrites the top summary rows, inserts four new percentile
# rows, drops four now-redundant percentile rows, and collapses the three
# multi-run "raw log line" rows down to their lead (summary) value.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1-3: top three summary values become "0M" (placeholder) ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# --- Step 4: row 4 "3" -> "386" ---
$t.Rows.Item(4).Cells.Item(1).Range.Text = "386"

# --- Step 5: insert four new rows right after row 4 (before the row that
# currently holds "0.00004"). Rows.Add(beforeRow) always inserts
# immediately above the *same* anchor row, so inserting in the reverse of
# the desired top-to-bottom order yields the correct final order. ---
$anchor = $t.Rows.Item(5)
$n1 = $t.Rows.Add($anchor)
$n1.Cells.Item(1).Range.Text = "0.00002"
$n2 = $t.Rows.Add($anchor)
$n2.Cells.Item(1).Range.Text = "0.00005"
$n3 = $t.Rows.Add($anchor)
$n3.Cells.Item(1).Range.Text = "0.00011"
$n4 = $t.Rows.Add($anchor)
$n4.Cells.Item(1).Range.Text = "0.00002"

# Rows 5 and 6 (originally "0.00004"/"0.00005") now sit at indices 9 and 10
# and keep their text. Rows 7 and 8 (originally "0.00005"/"0.00000", now at
# indices 11 and 12) get new values.
$t.Rows.Item(11).Cells.Item(1).Range.Text = "0.00011"
$t.Rows.Item(12).Cells.Item(1).Range.Text = "0.01773"

# --- Step 10: delete the four now-redundant rows that followed (originally
# "0.00004", "0.00005", "0.00005", "0.00014" — now at indices 13-16) ---
$t.Rows.Item(13).Delete()
$t.Rows.Item(13).Delete()
$t.Rows.Item(13).Delete()
$t.Rows.Item(13).Delete()

# --- Steps 12-14: collapse the three multi-run "raw log" rows (unaffected
# by the insert/delete above since they net to zero) down to a single run
# holding just the lead value. ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.98"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.02"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "83"
